$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenMap")
$tsMap = $wb.Worksheets.Item("timeslice map")

# --- New block headers (D1 / E1) used by the two new scenario-code blocks ---
$ws.Range("D1").Value = "vstacks_ts16~"
$ws.Range("E1").Value = "vstacks_t_annual~"

# --- Row 27-33: "ts-16" block (anchored on column D) ---
$ws.Range("A27").Formula = '=$D$1&TEXT(N27,"0000")'
for ($i = 28; $i -le 33; $i++) {
    $ws.Range("A$i").Formula = '=$D$1&TEXT(N' + $i + ',"0000")'
}

# --- Row 34-40: "ts-annual" block (anchored on column E) ---
$ws.Range("A34").Formula = '=$E$1&TEXT(N34,"0000")'
for ($i = 35; $i -le 40; $i++) {
    $ws.Range("A$i").Formula = '=$E$1&TEXT(N' + $i + ',"0000")'
}

# --- Common columns B, G, H, I, N, P for the 14 new rows (27-40) ---
for ($i = 27; $i -le 40; $i++) {
    $ws.Range("B$i").Formula = "=G$i"
    $ws.Range("G$i").Formula = "=H$i" + "&P$i"
    $src = $i - 7
    $ws.Range("H$i").Formula = "=H$src"
    $ws.Range("N$i").Formula = "=N$src"
}

# I / P columns: "ts-16" / "_16" for rows 27-33, "ts-annual" / "_ann" for rows 34-40
# (write all of column I first, then all of column P, matching the order the
# new shared-string table entries were authored in)
for ($i = 27; $i -le 33; $i++) {
    $ws.Range("I$i").Value = "ts-16"
}
for ($i = 34; $i -le 40; $i++) {
    $ws.Range("I$i").Value = "ts-annual"
}
for ($i = 27; $i -le 33; $i++) {
    $ws.Range("P$i").Value = "_16"
}
for ($i = 34; $i -le 40; $i++) {
    $ws.Range("P$i").Value = "_ann"
}

# --- Column A width widened slightly to fit the longer new codes ---
$ws.Columns.Item(1).ColumnWidth = 14.833333333333334

# --- Tab selection: ScenMap becomes the active / selected tab, "timeslice map" no longer is ---
$ws.Range("A33").Select()
$ws.Activate()

